# Applies the commit's row-reordering + "Förändrad" date bump to the
# "Avverkningsanmälningar" sheet.
#
# The rows in range 2..38 get permuted (several small blocks of rows were
# re-sorted in the source data), and every row's column C ("Förändrad")
# value changes from 46059 to 46060.
#
# Strategy: read the whole data block (rows 2-38) into memory first --
# columns A..R via .Value2 (raw values: plain text, numeric serials for
# dates instead of formatted DateTime strings, and the multi-line R column
# text, with blank cells coming back as plain $null) and columns S..Z via
# .Formula (the HYPERLINK formulas) -- then write back the permuted rows in
# one shot. Reading everything before writing anything avoids any chance of
# clobbering source data mid-way through the permutation. The two groups of
# columns are written back to two separate target ranges (A:R and S:Z) so
# that writing the all-formula array never touches the plain-value columns
# (and vice versa) -- writing a .Formula array that has blanks/nulls in
# some slots onto cells that also fall outside S:Z would otherwise wipe out
# the plain values just written there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 38
$nRows = $lastRow - $firstRow + 1   # 37

$valueRange = $ws.Range("A${firstRow}:R${lastRow}")      # columns 1..18
$formulaRange = $ws.Range("S${firstRow}:Z${lastRow}")    # columns 19..26

$valArr = $valueRange.Value2        # 1-based: [row 1..nRows, col 1..18]
$formulaArr = $formulaRange.Formula # 1-based: [row 1..nRows, col 1..8]

# Mapping of new row index (1-based, 1 == worksheet row 2) to the source row
# index (same numbering) it should take its content from.
$rowSource = @{
  1=1; 2=2; 3=3; 4=4; 5=5; 6=6; 7=7; 8=8; 9=9; 10=10;
  11=12; 12=11;
  13=13; 14=14; 15=15; 16=16; 17=17;
  18=20; 19=18; 20=19;
  21=22; 22=23; 23=21;
  24=24; 25=25;
  26=29; 27=26; 28=27; 29=28;
  30=30; 31=31;
  32=34; 33=32; 34=35; 35=36; 36=37; 37=33
}

$nValueCols = 18     # A..R
$nFormulaCols = 8    # S..Z

# Freshly created .NET arrays are 0-based, and must exactly match the target
# range's dimensions (nRows x nCols) or Excel will misalign the data.
$newValArr = New-Object 'object[,]' $nRows, $nValueCols
$newFormulaArr = New-Object 'object[,]' $nRows, $nFormulaCols

for ($newIdx = 1; $newIdx -le $nRows; $newIdx++) {
  $srcIdx = $rowSource[$newIdx]
  for ($col = 1; $col -le $nValueCols; $col++) {
    $v = $valArr[$srcIdx, $col]
    if ($null -eq $v) { $v = "" }
    $newValArr[$newIdx - 1, $col - 1] = $v
  }
  for ($col = 1; $col -le $nFormulaCols; $col++) {
    $f = $formulaArr[$srcIdx, $col]
    if ($null -eq $f) { $f = "" }
    $newFormulaArr[$newIdx - 1, $col - 1] = $f
  }
}

$valueRange.Value = $newValArr
$formulaRange.Formula = $newFormulaArr

# Column C ("Förändrad") becomes 46060 for every row in the block, regardless
# of the row permutation above (every single row had the same before/after
# value for this column).
$ws.Range("C${firstRow}:C${lastRow}").Value = 46060
